$d = $word.ActiveDocument

# Step 1: replace the first comma run (", ") with " and"
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Replacement.ClearFormatting()
$find1.Execute("LOS, to ensure", $true, $false, $false, $false, $false, $true, 1, $false, "LOS and to ensure", 2)

# Step 2: replace the trailing text run with just "."
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("care, and to ensure that we are legally allowed to use any predictors.", $true, $false, $false, $false, $false, $true, 1, $false, "care.", 2)
